# Apply cryptos.xlsx price/volume updates (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.889.81"
$ws.Range("E2").Value = "  -0.18%  "

# Row 3
$ws.Range("D3").Value = "2.279.86"
$ws.Range("E3").Value = "  -0.25%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "250.30"
$ws.Range("E5").Value = "  -0.77%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.632"
$ws.Range("E6").Value = "  +0.09%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "79.05"
$ws.Range("E7").Value = "  +8.37%  "

# Row 8
$ws.Range("E8").Value = "  +0.00%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.637"
$ws.Range("E9").Value = "  -3.60%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.35"
$ws.Range("E10").Value = "  +5.71%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0968"
$ws.Range("E11").Value = "  -1.26%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.36"
$ws.Range("E12").Value = "  -0.49%  "

# Row 13
$ws.Range("E13").Value = "  -1.25%  "

# Row 14
$ws.Range("D14").Value = "2.618.98"
$ws.Range("E14").Value = "  -0.27%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.18"
$ws.Range("E15").Value = "  +0.47%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.868"
$ws.Range("E16").Value = "  -2.40%  "

# Row 17
$ws.Range("D17").Value = "2.278.28"
$ws.Range("E17").Value = "  -0.32%  "

# Row 18
$ws.Range("D18").Value = "42.763.05"
$ws.Range("E18").Value = "  -0.34%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000100"
$ws.Range("E19").Value = "  -1.27%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.24"
$ws.Range("E20").Value = "  -2.42%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.24"
$ws.Range("E21").Value = "  -1.75%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "234.20"
$ws.Range("E22").Value = "  -0.70%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.18"
$ws.Range("E23").Value = "  +1.89%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.79"
$ws.Range("E24").Value = "  -2.61%  "

# Row 25
$ws.Range("E25").Value = "  +0.05%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.30"
$ws.Range("E26").Value = "  -3.40%  "

# Row 27
$ws.Range("E27").Value = "  -4.69%  "

# Row 28
$ws.Range("E28").Value = "  +1.86%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "168.90"
$ws.Range("E29").Value = "  +0.18%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.90"
$ws.Range("E30").Value = "  -1.22%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.69"
$ws.Range("E31").Value = "  +5.89%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0854"
$ws.Range("E32").Value = "  +4.87%  "

# Row 33
$ws.Range("E33").Value = "  -5.14%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "30.55"
$ws.Range("E34").Value = "  -2.34%  "

# Row 35
$ws.Range("E35").Value = "  +1.25%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.58"
$ws.Range("E36").Value = "  -4.60%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.78"
$ws.Range("E37").Value = "  -0.51%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0304"
$ws.Range("E38").Value = "  -2.10%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.58"
$ws.Range("E39").Value = "  +2.53%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.27"
$ws.Range("E40").Value = "  -2.88%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.94"
$ws.Range("E41").Value = "  -1.77%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "115.07"
$ws.Range("E42").Value = "  +17.92%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.211"
$ws.Range("E43").Value = "  -1.76%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "61.52"
$ws.Range("E44").Value = "  -0.49%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.92"
$ws.Range("E45").Value = "  -3.32%  "

# Row 46
$ws.Range("B46").Value = "FTXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.71"
$ws.Range("E46").Value = "  -5.72%  "

# Row 47
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.102"
$ws.Range("E47").Value = "  -2.52%  "

# Row 48
$ws.Range("E48").Value = "  +0.01%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.16"
$ws.Range("E49").Value = "  -3.53%  "

# Row 50
$ws.Range("E50").Value = "  -2.17%  "

# Row 51
$ws.Range("E51").Value = "  -0.96%  "
